# Add a new "Concept_3" column (E) to the vehicle characteristics sheet.
# Column E mirrors column D (Concept_2) for every parameter, except the
# Cx value (row 19): Concept_2 (D19) changes from 2 to 3, and the new
# Concept_3 (E19) is 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("E2").Value = "Concept_3"

# --- Plain numeric parameter values (rows 3-18), copied from column D ---
$ws.Range("E3").Value = 1.236
$ws.Range("E4").Value = 1.165
$ws.Range("E5").Value = 1.575
$ws.Range("E6").Value = 0.78749999999999998
$ws.Range("E7").Value = 0.52
$ws.Range("E8").Value = 0.01
$ws.Range("E9").Value = 1.5
$ws.Range("E10").Value = 0.77110000000000001
$ws.Range("E11").Value = 0.00005
$ws.Range("E12").Value = 0.005
$ws.Range("E13").Value = 0.03
$ws.Range("E14").Value = 0.313
$ws.Range("E15").Value = 70
$ws.Range("E16").Value = 215
$ws.Range("E17").Value = 0.42
$ws.Range("E18").Value = 0.31

# --- Cx (row 19): Concept_2 becomes 3, new Concept_3 is 4 ---
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 4

# --- Remaining SUSPENSION block (rows 20-24) ---
$ws.Range("E20").Value = 2.13
$ws.Range("E21").Value = 0.5
$ws.Range("E22").Value = 1.1399999999999999
$ws.Range("E23").Value = 0.4
$ws.Range("E24").Value = 15

# --- GENERAL block (row 25 engine type, rows 26-38 numeric) ---
$ws.Range("E25").Value = "CBR600RR"
$ws.Range("E26").Value = 0.9
$ws.Range("E27").Value = 9000
$ws.Range("E28").Value = 13500
$ws.Range("E29").Value = 0.5
$ws.Range("E30").Value = 1111
$ws.Range("E31").Value = 0.005
$ws.Range("E32").Value = 0.1
$ws.Range("E33").Value = 1111
$ws.Range("E34").Value = 0.3125
$ws.Range("E35").Value = 6
$ws.Range("E36").Value = 0.47370000000000001
$ws.Range("E37").Value = 1111
$ws.Range("E38").Value = 1111

# --- Formula rows (39-40), recreated with column E references ---
$ws.Range("E40").Formula = "=E15+E16+E24"
$ws.Range("E39").Formula = "=(E17*E15+E18*E16+E23*E24)/E40"

# --- Row 41 ---
$ws.Range("E41").Value = 0.1

# Move the active selection, matching the saved view state
$ws.Range("F19").Select()

$wb.Save()
